$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.758.13'
$ws.Range("E2").Value = '  +2.24%  '
$ws.Range("D3").Value = '2.942.54'
$ws.Range("E3").Value = '  +0.42%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = "'593.14"
$ws.Range("E5").Value = '  -0.76%  '
$ws.Range("D6").Value = "'147.25"
$ws.Range("E6").Value = '  +1.32%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '2.940.93'
$ws.Range("E8").Value = '  +0.40%  '
$ws.Range("D9").Value = "'0.506"
$ws.Range("E9").Value = '  +0.80%  '
$ws.Range("E10").Value = '  +3.81%  '
$ws.Range("E11").Value = '  +5.51%  '
$ws.Range("D12").Value = "'0.441"
$ws.Range("E12").Value = '  +0.56%  '
$ws.Range("E13").Value = '  +4.66%  '
$ws.Range("D14").Value = "'32.64"
$ws.Range("E14").Value = '  -2.61%  '
$ws.Range("E15").Value = '  -0.71%  '
$ws.Range("D16").Value = '3.430.05'
$ws.Range("E16").Value = '  +0.40%  '
$ws.Range("D17").Value = '62.707.09'
$ws.Range("E17").Value = '  +2.12%  '
$ws.Range("D18").Value = "'6.68"
$ws.Range("E18").Value = '  -0.15%  '
$ws.Range("D19").Value = '2.942.96'
$ws.Range("E19").Value = '  +0.41%  '
$ws.Range("D20").Value = "'439.56"
$ws.Range("E20").Value = '  +1.87%  '
$ws.Range("D21").Value = "'13.36"
$ws.Range("E21").Value = '  -0.72%  '
$ws.Range("D22").Value = "'0.665"
$ws.Range("E22").Value = '  -1.68%  '
$ws.Range("D23").Value = "'7.03"
$ws.Range("E23").Value = '  -0.76%  '
$ws.Range("D24").Value = "'80.86"
$ws.Range("E24").Value = '  -1.30%  '
$ws.Range("D25").Value = "'11.12"
$ws.Range("E25").Value = '  +2.24%  '
$ws.Range("D26").Value = "'2.13"
$ws.Range("E26").Value = '  -2.78%  '
$ws.Range("D27").Value = "'11.74"
$ws.Range("E27").Value = '  -0.28%  '
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("D29").Value = "'2.24"
$ws.Range("E29").Value = '  +0.63%  '
$ws.Range("D30").Value = "'7.17"
$ws.Range("E30").Value = '  +3.59%  '
$ws.Range("D31").Value = "'2.61"
$ws.Range("E31").Value = '  -0.18%  '
$ws.Range("E32").Value = '  +14.85%  '
$ws.Range("D33").Value = "'0.109"
$ws.Range("E33").Value = '  -0.89%  '
$ws.Range("D34").Value = "'26.32"
$ws.Range("E34").Value = '  -1.11%  '
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = '  -0.16%  '
$ws.Range("D36").Value = "'0.991"
$ws.Range("E36").Value = '  -2.10%  '
$ws.Range("D37").Value = "'3.09"
$ws.Range("E37").Value = '  +3.66%  '
$ws.Range("D38").Value = "'5.56"
$ws.Range("E38").Value = '  -1.14%  '
$ws.Range("D39").Value = "'49.66"
$ws.Range("E39").Value = '  -0.67%  '
$ws.Range("D40").Value = "'2.02"
$ws.Range("E40").Value = '  +0.98%  '
$ws.Range("D41").Value = "'8.46"
$ws.Range("E41").Value = '  -1.35%  '
$ws.Range("E42").Value = '  -4.48%  '
$ws.Range("D43").Value = "'0.280"
$ws.Range("E43").Value = '  -0.84%  '
$ws.Range("D44").Value = "'39.01"
$ws.Range("E44").Value = '  -7.85%  '
$ws.Range("D45").Value = '2.695.82'
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").Value = "'134.51"
$ws.Range("E46").Value = '  +0.22%  '
$ws.Range("D47").Value = "'361.50"
$ws.Range("E47").Value = '  +0.37%  '
$ws.Range("D48").Value = "'0.0336"
$ws.Range("E48").Value = '  -2.85%  '
$ws.Range("E50").Value = '  -0.77%  '
$ws.Range("D51").Value = "'22.66"
$ws.Range("E51").Value = '  -4.12%  '
